# Update benchmarking stats for PEPMatch (row 2) and the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PEPMatch benchmark values (row 2)
$ws.Range("B2").Value = 13.381
$ws.Range("D2").Value = 20.513
$ws.Range("E2").Value = 33.894

# Normalize the formatting on row 2 back to the default style
$ws.Range("B2:G2").Font.Name = "Calibri"

# Row 2 height shrinks slightly after the refresh
$ws.Rows(2).RowHeight = 13.8

# Move the active selection to G3
$ws.Range("G3").Select()
